$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.348.10"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.879.21"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.93"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.682"
$ws.Range("E6").Value = "  -2.33%  "
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.62"
$ws.Range("E8").Value = "  +5.03%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.32"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0744"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "13.50"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "2.152.32"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.764"
$ws.Range("E15").Value = "  +4.66%  "
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").Value = "1.900.13"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "35.352.09"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.66"
$ws.Range("E19").Value = "  -0.94%  "
$ws.Range("D20").Value = "0.0₃0825"
$ws.Range("E20").Value = "  -1.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.49"
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.82"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.07"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").Value = "  +10.60%  "
$ws.Range("E25").Value = "  -0.54%  "
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.64"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.66"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("E30").Value = "  -1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.29"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0592"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("E34").Value = "  -0.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.79"
$ws.Range("E35").Value = "  -10.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.45"
$ws.Range("E36").Value = "  -10.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.856"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0729"
$ws.Range("E39").Value = "  +9.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.40"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.81"
$ws.Range("E42").Value = "  -2.11%  "
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "1.307.18"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0801"
$ws.Range("E46").Value = "  +4.13%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.88"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("E50").Value = "  -4.44%  "
$ws.Range("D51").Value = "2.055.82"
$ws.Range("E51").Value = "  -1.35%  "
